# Apply MCDA_baseline.xlsx updates: 'Fixed T9, ran analysis, new graphs'
# Updates Alternative B (E) and Alternative C (F) scores on 'Score' sheet
# and corresponding ranks (D/E/F) on 'Rank' sheet, per the re-run analysis.

$wb = $excel.ActiveWorkbook
$wsScore = $wb.Worksheets.Item("Score")
$wsRank = $wb.Worksheets.Item("Rank")

# --- Score sheet updates (Alternative B / Alternative C columns) ---
$wsScore.Range("E2").Value = 0.2227685477536643
$wsScore.Range("F2").Value = 0.7772314522463357
$wsScore.Range("E7").Value = 0.8643220542436826
$wsScore.Range("F7").Value = 0.1446208252063952
$wsScore.Range("E8").Value = 0.5434327177894198
$wsScore.Range("F8").Value = 0.4733938948918839
$wsScore.Range("E9").Value = 0.5738460381881468
$wsScore.Range("F9").Value = 0.4509750594913713
$wsScore.Range("E10").Value = 0.8598086312293285
$wsScore.Range("F10").Value = 0.1494026779803342
$wsScore.Range("E11").Value = 0.5188625936341134
$wsScore.Range("F11").Value = 0.5023898737537218
$wsScore.Range("E12").Value = 0.9002138881024789
$wsScore.Range("F12").Value = 0.09978611189752094
$wsScore.Range("E13").Value = 0.7244050160684078
$wsScore.Range("F13").Value = 0.2755949839315923
$wsScore.Range("E14").Value = 0.5738460381881468
$wsScore.Range("F14").Value = 0.4509750594913713
$wsScore.Range("E15").Value = 0.8598086312293285
$wsScore.Range("F15").Value = 0.1494026779803342
$wsScore.Range("E16").Value = 0.8788407057929699
$wsScore.Range("F16").Value = 0.1296768366373445
$wsScore.Range("E17").Value = 0.6942695365180909
$wsScore.Range("F17").Value = 0.3227899650669473
$wsScore.Range("E28").Value = 0.4233067709684146
$wsScore.Range("F28").Value = 0.5884870768418804
$wsScore.Range("E29").Value = 0.4079095200044974
$wsScore.Range("F29").Value = 0.6008026351603301
$wsScore.Range("E30").Value = 0.8187598662855726
$wsScore.Range("F30").Value = 0.1812401337144274
$wsScore.Range("E31").Value = 0.5724224368103152
$wsScore.Range("F31").Value = 0.4275775631896847
$wsScore.Range("E32").Value = 0.6195447123585156
$wsScore.Range("F32").Value = 0.4015400151485282
$wsScore.Range("E36").Value = 0.5862029576185354
$wsScore.Range("F36").Value = 0.4477419125491551
$wsScore.Range("E40").Value = 0.9474680615893214
$wsScore.Range("F40").Value = 0.05253193841067868
$wsScore.Range("E44").Value = 0.8390467454170564
$wsScore.Range("F44").Value = 0.1609532545829435

# --- Rank sheet updates ---
$wsRank.Range("D2").Value = 2
$wsRank.Range("E2").Value = 3
$wsRank.Range("F2").Value = 1
$wsRank.Range("D11").Value = 3
$wsRank.Range("F11").Value = 2
$wsRank.Range("E28").Value = 2
$wsRank.Range("F28").Value = 1
$wsRank.Range("E29").Value = 3
$wsRank.Range("F29").Value = 1
